$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Groups")

# Make room for a new "Stock" column by shifting the existing header
# cells (Year ... First Detail Event) one column to the right, from I2
# down to A2, so that column A is freed up without touching the sheet's
# <cols> column-width definitions.
for ($col = 9; $col -ge 1; $col--) {
    $src = $ws.Cells.Item(2, $col)
    $dst = $ws.Cells.Item(2, $col + 1)
    $dst.Value = $src.Value()
}

# J2 is a brand new cell - give it the same header formatting as the
# rest of the row before we overwrite A2.
$ws.Cells.Item(2, 9).Copy()
$ws.Cells.Item(2, 10).PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Populate the freed-up A2 cell with the new "Stock" column header.
$ws.Cells.Item(2, 1).Value = "Stock"

# Match the recorded selection in the sheet.
$ws.Activate()
$ws.Range("B5").Select() | Out-Null
